$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 644, shifting existing rows 644:745 down to 645:746
$ws.Rows("644:644").Insert()

# Populate the newly inserted row 644 with its data
$ws.Range("A644").Value = 4
$ws.Range("B644").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C644").Value = "Los Lagos"
$ws.Range("D644").Value = 45034
$ws.Range("E644").Value = 10
$ws.Range("F644").Value = "Fruta"
$ws.Range("G644").Value = 100102
$ws.Range("H644").Value = "Cítricos"
$ws.Range("I644").Value = 100102005
$ws.Range("J644").Value = "Naranja"
$ws.Range("K644").Value = "Valencia"
$ws.Range("L644").Value = "Primera"
$ws.Range("M644").Value = 600
$ws.Range("N644").Value = 19000
$ws.Range("O644").Value = 20000
$ws.Range("P644").Value = 19500
$ws.Range("Q644").Value = '$/caja 15 kilos empedrada'
$ws.Range("R644").Value = "Región de O'Higgins"
$ws.Range("S644").Value = 1300
$ws.Range("T644").Value = 15
